# Agregar datos al Excel para las prueba TC0012
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataPrueba")

$ws.Range("A13").Value = "TC0012_Filtrado_Trenes_MAD_VLC_Precio_Equipaje_Escala_Salida_Aere_Esta"
$ws.Range("B13").Value = "Madrid"
$ws.Range("C13").Value = "Valencia"
$ws.Range("D13").Value = 16
$ws.Range("E13").Value = 20

$ws.Range("A13").Select()
